$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.922.15'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.884.83'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.18'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4591'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3891'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07849'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9863'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.84'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.96'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.38%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.76%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06939'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.20'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.48%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009968'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.915.36'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.276'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.096.31'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.089'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.03'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.28'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.990'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.96%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.45'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09326'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9046'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.279'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.327'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.267'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.203'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05769'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02072'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.20%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.636'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5667'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1768'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.711'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.259'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.90'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5360'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07037'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.850'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.38%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Quant'

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.77'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'MXToken'

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.532'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.070'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.93%  '
